$wb = $excel.ActiveWorkbook

# --- "location" sheet: seed it with x/y coordinates for each zone, and
# rename the generic "id" header to the explicit "location_id" ---
$ws = $wb.Worksheets.Item("location")
$ws.Activate()

$ws.Range("C1").Value = "x_coordinate"
$ws.Range("D1").Value = "y_coordinate"
$ws.Range("A1").Value = "location_id"

$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 3

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1

# auto-fit the columns that now hold new/renamed data
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# leave the selection on the first data row, same as the saved workbook
$null = $ws.Range("A2").Select()
